# Capitalizar nombres de puntos para extraerlos desde acá
# Row 2 of Sheet1 contains street/point names used as column headers.
# Several of them were written in lowercase ("avenida jiménez", "calle 32", ...)
# and need to be capitalized ("Avenida Jiménez", "Calle 32", ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "H2"  = "Avenida Jiménez"
    "K2"  = "Calle 32"
    "L2"  = "Calle 37"
    "O2"  = "Calle 53"
    "P2"  = "Calle 56"
    "Q2"  = "Calle 60"
    "R2"  = "Calle 67"
    "S2"  = "Calle 70"
    "T2"  = "Calle 72"
    "U2"  = "Calle 74"
    "V2"  = "Calle 76"
    "W2"  = "Calle 84"
    "X2"  = "Calle 85"
    "Y2"  = "Calle 92"
    "Z2"  = "Calle 94"
    "AA2" = "Calle 100"
    "AB2" = "Calle 106"
    "AC2" = "Calle 116"
    "AD2" = "Calle 127"
    "AE2" = "Calle 134"
    "AF2" = "Calle 140"
    "AG2" = "Calle 147"
    "AH2" = "Calle 151"
    "AI2" = "Calle 153"
    "AJ2" = "Calle 160"
    "AK2" = "Calle 164"
    "AL2" = "Calle 165"
    "AM2" = "Calle 170"
    "AN2" = "Calle 175"
    "AO2" = "Calle 180"
    "AP2" = "Calle 183"
    "AQ2" = "Calle 189"
    "AR2" = "Calle 192"
    "AS2" = "Calle 193"
    "AT2" = "Calle 200"
    "AU2" = "Calle 220"
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
